# Auto-update Cloudflare export
#
# The live Cloudflare DNS export changed between the two snapshots:
#   - the record "arena.repolizer.com" (row 109) was removed
#   - the record "stats.repolizer.com" (row 111, a CNAME) was removed
#   - a new record "www.repolizer.com" (a CNAME pointing at the same
#     cfargotunnel.com target as ssh.repolizer.com) was added right
#     after "ssh.repolizer.com"
# Every other row is unchanged and simply shifts up to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "arena.repolizer.com" record (row 109).
$ws.Rows("109:109").Delete()

# After the delete above, the old "stats.repolizer.com" CNAME record has
# shifted up into row 110 - remove it too.
$ws.Rows("110:110").Delete()

# Insert a fresh blank row to hold the new "www.repolizer.com" record,
# placed right after "ssh.repolizer.com" (now row 109).
$ws.Rows("110:110").Insert()

$ws.Range("A110").Value = "bfd43018c9b1cd41cc29efc4c62c6d8c"
$ws.Range("B110").Value = "www.repolizer.com"
$ws.Range("C110").Value = "CNAME"
$ws.Range("D110").Value = "a7a86777-918b-4a27-9ad1-ab108649b92d.cfargotunnel.com"
$ws.Range("E110").Value = $true
$ws.Range("F110").Value = $true
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = "{}"
$ws.Range("I110").Value = "{}"
$ws.Range("J110").Value = ""
$ws.Range("K110").Value = "[]"
$ws.Range("L110").Value = "2025-04-23T08:30:24.446765Z"
$ws.Range("M110").Value = "2025-04-23T08:34:01.380912Z"
$ws.Range("N110").Value = ""
$ws.Range("O110").Value = ""
